$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap name order: "First Last" -> "Last First"
$ws.Range("F3").Value = "Renzi Marco"
$ws.Range("G3").Value = "Testa Filippo"

# Update the active selection (cosmetic change from the original edit)
$ws.Range("B7").Select()
